$wb = $excel.ActiveWorkbook

# --- Sheet: Productdata ---
$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("E2").Value = 11.6942496
$ws.Range("E3").Value = 2.1712096
$ws.Range("E5").Value = 1.421352
$ws.Range("E6").Value = 0.9547008
$ws.Range("E7").Value = 0.2759328
$ws.Range("E8").Value = 0.09377280000000002
$ws.Range("E9").Value = 0.8618752000000001
$ws.Range("C10").Value = 601
$ws.Range("E10").Value = 0.5204736
$ws.Range("C11").Value = 1803
$ws.Range("E11").Value = 0.7561680000000001
$ws.Range("C12").Value = 423
$ws.Range("E12").Value = 1.2251888
$ws.Range("C13").Value = 3237
$ws.Range("E13").Value = 12.57585599999999
$ws.Range("C14").Value = 1190
$ws.Range("E14").Value = 4.844571200000001
$ws.Range("C15").Value = 303
$ws.Range("C16").Value = 417
$ws.Range("E16").Value = 0.8744224
$ws.Range("C17").Value = 597
$ws.Range("E17").Value = 1.392768
$ws.Range("C18").Value = 170
$ws.Range("E18").Value = 0.432432
$ws.Range("C19").Value = 59
$ws.Range("E19").Value = 0.1297664
$ws.Range("E20").Value = 62.875936
$ws.Range("E21").Value = 66.74890239999999
$ws.Range("E22").Value = 82.6232832
$ws.Range("E23").Value = 254.2291648

# --- Sheet: ForecastedAverageDemand ---
$ws = $wb.Worksheets.Item("ForecastedAverageDemand")
$ws.Range("C2").Value = 596
$ws.Range("D2").Value = 150
$ws.Range("F2").Value = 295
$ws.Range("G2").Value = 83
$ws.Range("I2").Value = 422
$ws.Range("J2").Value = 301
$ws.Range("K2").Value = 899
$ws.Range("L2").Value = 209
$ws.Range("C3").Value = 602
$ws.Range("D3").Value = 145
$ws.Range("F3").Value = 289
$ws.Range("G3").Value = 87
$ws.Range("H3").Value = 31
$ws.Range("I3").Value = 418
$ws.Range("K3").Value = 904
$ws.Range("L3").Value = 214
$ws.Range("C4").Value = 604
$ws.Range("D4").Value = 154
$ws.Range("F4").Value = 296
$ws.Range("G4").Value = 89
$ws.Range("H4").Value = 26
$ws.Range("I4").Value = 418
$ws.Range("J4").Value = 299
$ws.Range("K4").Value = 899
$ws.Range("L4").Value = 202
$ws.Range("C5").Value = 586
$ws.Range("D5").Value = 149
$ws.Range("F5").Value = 301
$ws.Range("G5").Value = 81
$ws.Range("H5").Value = 33
$ws.Range("I5").Value = 421
$ws.Range("J5").Value = 302
$ws.Range("K5").Value = 898
$ws.Range("L5").Value = 215
$ws.Range("C6").Value = 601
$ws.Range("D6").Value = 148
$ws.Range("F6").Value = 307
$ws.Range("G6").Value = 89
$ws.Range("H6").Value = 29
$ws.Range("I6").Value = 417
$ws.Range("J6").Value = 304
$ws.Range("K6").Value = 901
$ws.Range("L6").Value = 199

# --- Sheet: ForcastedStandardDeviation ---
$ws = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws.Range("C2").Value = 74.5
$ws.Range("D2").Value = 18.75
$ws.Range("F2").Value = 36.875
$ws.Range("G2").Value = 10.375
$ws.Range("I2").Value = 52.75
$ws.Range("J2").Value = 37.625
$ws.Range("K2").Value = 112.375
$ws.Range("L2").Value = 26.125
$ws.Range("C3").Value = 112.875
$ws.Range("D3").Value = 27.1875
$ws.Range("F3").Value = 54.1875
$ws.Range("G3").Value = 16.3125
$ws.Range("H3").Value = 5.8125
$ws.Range("I3").Value = 78.375
$ws.Range("K3").Value = 169.5
$ws.Range("L3").Value = 40.125
$ws.Range("C4").Value = 132.125
$ws.Range("D4").Value = 33.6875
$ws.Range("F4").Value = 64.75
$ws.Range("G4").Value = 19.46875
$ws.Range("H4").Value = 5.6875
$ws.Range("I4").Value = 91.4375
$ws.Range("J4").Value = 65.40625
$ws.Range("K4").Value = 196.65625
$ws.Range("L4").Value = 44.1875
$ws.Range("C5").Value = 137.34375
$ws.Range("D5").Value = 34.921875
$ws.Range("F5").Value = 70.546875
$ws.Range("G5").Value = 18.984375
$ws.Range("H5").Value = 7.734375
$ws.Range("I5").Value = 98.671875
$ws.Range("J5").Value = 70.78125
$ws.Range("K5").Value = 210.46875
$ws.Range("L5").Value = 50.390625
$ws.Range("C6").Value = 145.5546875
$ws.Range("D6").Value = 35.84375
$ws.Range("F6").Value = 74.3515625
$ws.Range("G6").Value = 21.5546875
$ws.Range("H6").Value = 7.0234375
$ws.Range("I6").Value = 100.9921875
$ws.Range("J6").Value = 73.625
$ws.Range("K6").Value = 218.2109375
$ws.Range("L6").Value = 48.1953125

# --- Sheet: Capacity ---
$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B2").Value = 50428
$ws.Range("B3").Value = 522970
$ws.Range("B4").Value = 522970
$ws.Range("B5").Value = 4930859.999999999
